$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email value in D2
$ws.Range("D2").Value = "pooja168@givmail.com"

# Move the active selection to D7 (matches the captured UI state after the edit)
$ws.Range("D7").Select()
